$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "Save" column - copy G1's formatting (border/font/alignment)
# onto H1 so it reuses the existing header style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New "Save" column values for the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
